# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 454
$wsOff.Range("C2").Value = 347
$wsOff.Range("D2").Value = 138
$wsOff.Range("E2").Value = 68

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 494
$wsDef.Range("C2").Value = 353
$wsDef.Range("D2").Value = 103
$wsDef.Range("E2").Value = 43
$wsDef.Range("F2").Value = 6
$wsDef.Range("G2").Value = 10
